$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 16
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()

# Row 43
$ws.Range("H43").Value = 3849153.5
$ws.Range("I43").Value = 15384615
$ws.Range("K43").Value = 15384615
$ws.Range("M43").Value = -15384546

# Row 132
$ws.Range("H132").Value = 135922.7
$ws.Range("I132").Value = 183421.62
$ws.Range("K132").Value = 550264.86
$ws.Range("M132").Value = -547734.86

# Row 137
$ws.Range("H137").Value = 504467.7
$ws.Range("I137").Value = 1251996.4
$ws.Range("J137").Value = 6115.25
$ws.Range("K137").Value = 3755989.2
$ws.Range("L137").Value = 18345.75
$ws.Range("M137").Value = -3753439.2
$ws.Range("N137").Value = -23445.75

# Row 138
$ws.Range("H138").Value = 2558.899
$ws.Range("I138").Value = 2034.0555
$ws.Range("J138").Value = 2691.9578
$ws.Range("K138").Value = 6102.166499999999
$ws.Range("L138").Value = 8075.8734
$ws.Range("M138").Value = -962.1664999999994
$ws.Range("N138").Value = -18355.8734

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 14681.349
$ws.Range("I32").Value = 14030.019
$ws.Range("J32").Value = 17938
$ws.Range("K32").Value = 14030.019
$ws.Range("L32").Value = 17938
$ws.Range("M32").Value = -13743.019
$ws.Range("N32").Value = -18512

# Row 97
$ws.Range("H97").Value = 764.6429000000001
$ws.Range("I97").Value = 746.53845
$ws.Range("K97").Value = 746.53845
$ws.Range("M97").Value = -250.53845

# Row 110
$ws.Range("H110").Value = 2046503.2
$ws.Range("I110").Value = 2553504.5
$ws.Range("K110").Value = 2553504.5
$ws.Range("M110").Value = -2551459.5

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 1864.3334
$ws.Range("I105").Value = 1864.3334
$ws.Range("K105").Value = 1864.3334
$ws.Range("M105").Value = -117.3334

$ws = $wb.Worksheets.Item("CRP")
# Row 11
$ws.Range("H11").Value = 121293
$ws.Range("I11").Value = 500
$ws.Range("J11").Value = 151491.25
$ws.Range("K11").Value = 500
$ws.Range("L11").Value = 151491.25
$ws.Range("M11").Value = -360
$ws.Range("N11").Value = -151771.25

# Row 58
$ws.Range("H58").Value = 2787.5
$ws.Range("I58").Value = 1519.3182
$ws.Range("K58").Value = 1519.3182
$ws.Range("M58").Value = -1316.3182

# Row 134
$ws.Range("H134").Value = 1916.4517
$ws.Range("I134").Value = 1927.56
$ws.Range("K134").Value = 5782.68
$ws.Range("M134").Value = -3247.68

# Row 136
$ws.Range("H136").Value = 2787.5
$ws.Range("I136").Value = 1519.3182
$ws.Range("K136").Value = 4557.9546
$ws.Range("M136").Value = -2007.9546

$ws = $wb.Worksheets.Item("CUL")
# Row 37
$ws.Range("H37").Value = 68234.5
$ws.Range("J37").Value = 68234.5
$ws.Range("L37").Value = 204703.5
$ws.Range("N37").Value = -204927.5

# Row 64
$ws.Range("H64").Value = 19610
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 19610
$ws.Range("K64").Value = 0
$ws.Range("L64").ClearContents()
$ws.Range("M64").Value = 58830
$ws.Range("N64").Value = -59370

# Row 67
$ws.Range("H67").Value = 19610
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 19610
$ws.Range("K67").Value = 0
$ws.Range("L67").ClearContents()
$ws.Range("M67").Value = 58830
$ws.Range("N67").Value = -60702

# Row 69
$ws.Range("H69").Value = 7000
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()

# Row 72
$ws.Range("H72").Value = 7000
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()

# Row 81
$ws.Range("H81").Value = 5082.875

# Row 84
$ws.Range("H84").Value = 5082.875

# Row 98
$ws.Range("H98").Value = 1308.5625
$ws.Range("J98").Value = 1359.7858
$ws.Range("L98").Value = 4079.3574
$ws.Range("N98").Value = -7075.357400000001

# Row 103
$ws.Range("H103").Value = 575.2
$ws.Range("I103").Value = 259.6
$ws.Range("K103").Value = 778.8000000000001
$ws.Range("M103").Value = 100.1999999999999

# Row 131
$ws.Range("H131").Value = 12312728
$ws.Range("I131").Value = 55555852
$ws.Range("J131").Value = 9429853
$ws.Range("K131").Value = 166667556
$ws.Range("L131").Value = 28289559
$ws.Range("M131").Value = -166662516
$ws.Range("N131").Value = -28299639

# Row 132
$ws.Range("H132").Value = 3689.8635
$ws.Range("I132").Value = 1007.875
$ws.Range("J132").Value = 5222.4287
$ws.Range("K132").Value = 9070.875
$ws.Range("L132").Value = 47001.85830000001
$ws.Range("M132").Value = -6540.875
$ws.Range("N132").Value = -52061.85830000001

# Row 137
$ws.Range("H137").Value = 170002620
$ws.Range("I137").Value = 187503280
$ws.Range("J137").Value = 100000000
$ws.Range("K137").Value = 562509840
$ws.Range("L137").Value = 300000000
$ws.Range("M137").Value = -562504740
$ws.Range("N137").Value = -300010200

$ws = $wb.Worksheets.Item("GSM")
# Row 15
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").ClearContents()
$ws.Range("M15").ClearContents()
$ws.Range("N15").Value = 0

# Row 29
$ws.Range("H29").Value = 15000
$ws.Range("J29").Value = 15000
$ws.Range("L29").Value = 15000
$ws.Range("N29").Value = -15580

# Row 81
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").ClearContents()
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = 0

# Row 84
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").ClearContents()
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = 0

$ws = $wb.Worksheets.Item("LTW")
# Row 32
$ws.Range("H32").Value = 11400.333
$ws.Range("I32").Value = 11400.333
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 11400.333
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -11083.333

$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 7845.2856
$ws.Range("I136").Value = 2919
$ws.Range("K136").Value = 8757
$ws.Range("M136").Value = -6207

# Row 137
$ws.Range("H137").Value = 109889
$ws.Range("J137").Value = 109889
$ws.Range("L137").Value = 109889
$ws.Range("N137").Value = -120089
